# 增加 LocationComponent 组件 未完成
# Update the AllServer row's external (outer) address from the old
# 127.0.0.1:10001 placeholder to 127.0.0.1:10025, and move the active
# selection to L5 as left by the author when they saved the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "127.0.0.1:10025"

$ws.Range("L5").Select()
